$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Akram%M Zaaqoq%coreGivesNoEmail%1,          Akshay%Kohli%coreGivesNoEmail%1,          Alexandra%Pratt%coreGivesNoEmail%1,          Emil%Oweis%coreGivesNoEmail%1,          Eric%Kriner%coreGivesNoEmail%1,          Ezequiel%Molina%coreGivesNoEmail%1,          Fatima%Hayat%coreGivesNoEmail%1,          Karthik%Vedantam%coreGivesNoEmail%1,          Makdisi%G%coreGivesNoEmail%1,          Muhammad%D Hashmi%coreGivesNoEmail%1,          Muhtadi%Alnababteh%coreGivesNoEmail%1,          Rajus%Chopra%coreGivesNoEmail%1,          World%Health Organization%coreGivesNoEmail%1,          Yang%X%coreGivesNoEmail%1]'
$ws.Range("E3").Value = '[Waqas%Akhtar%NULL%1,           Olusegun%Olusanya%NULL%2,           Olusegun%Olusanya%NULL%0,           Marta Montero%Baladia%NULL%1,           Harriet%Young%NULL%1,           Sachin%Shah%s.shah4@nhs.net%1]'
$ws.Range("E4").Value = '[Ryan P%Barbaro%NULL%0,           Graeme%MacLaren%NULL%1,           Philip S%Boonstra%NULL%1,           Theodore J%Iwashyna%NULL%1,           Arthur S%Slutsky%NULL%1,           Eddy%Fan%NULL%1,           Robert H%Bartlett%NULL%1,           Joseph E%Tonna%NULL%1,           Robert%Hyslop%NULL%1,           Jeffrey J%Fanning%NULL%1,           Peter T%Rycus%NULL%1,           Steve J%Hyer%NULL%1,           Marc M%Anders%NULL%1,           Cara L%Agerstrand%NULL%1,           Katarzyna%Hryniewicz%NULL%1,           Rodrigo%Diaz%NULL%1,           Roberto%Lorusso%NULL%1,           Alain%Combes%NULL%1,           Daniel%Brodie%NULL%0,           Peta%Alexander%NULL%1,           Nicholas%Barrett%NULL%0,           Jan%Bělohlávek%NULL%1,           Dale%Fisher%NULL%1,           John%Fraser%NULL%1,           Ali Ait%Hssain%NULL%0,           Jae Sung%Jung%NULL%1,           Michael%McMullan%NULL%1,           Yatin%Mehta%NULL%1,           Mark T.%Ogino%NULL%1,           Matthew L.%Paden%NULL%1,           Kiran%Shekar%NULL%1,           Christine%Stead%NULL%1,           Yasir%Abu-Omar%NULL%1,           Vanni%Agnoletti%NULL%1,           Anzila%Akbar%NULL%1,           Huda%Alfoudri%NULL%1,           Carlos%Alviar%NULL%1,           Vladimir%Aronsky%NULL%1,           Erin%August%NULL%1,           Georg%Auzinger%NULL%1,           Hilda%Aveja%NULL%1,           Rhonda%Bakken%NULL%1,           Joan%Balcells%NULL%1,           Sripal%Bangalore%NULL%1,           Bernard W.%Barnes%NULL%1,           Alaiza%Bautista%NULL%1,           Lorraine L.%Bellows%NULL%1,           Felipe%Beltran%NULL%1,           Peyman%Benharash%NULL%1,           Marco%Benni%NULL%1,           Jennifer%Berg%NULL%1,           Pietro%Bertini%NULL%1,           Pablo%Blanco-Schweizer%NULL%1,           Melissa%Brunsvold%NULL%1,           Jenny%Budd%NULL%1,           Debra%Camp%NULL%1,           Mark%Caridi-Scheible%NULL%0,           Edmund%Carton%NULL%1,           Elena%Casanova-Ghosh%NULL%1,           Anthony%Castleberry%NULL%1,           Christopher T.%Chipongian%NULL%1,           Chang Woo%Choi%NULL%1,           Alessandro%Circelli%NULL%1,           Elliott%Cohen%NULL%1,           Michael%Collins%NULL%1,           Scott%Copus%NULL%1,           Jill%Coy%NULL%1,           Brandon%Crist%NULL%1,           Leonora%Cruz%NULL%1,           Mirosław%Czuczwar%NULL%1,           Mani%Daneshmand%NULL%1,           Daniel%Davis II%NULL%1,           Kim%De la Cruz%NULL%1,           Cyndie%Devers%NULL%1,           Toni%Duculan%NULL%1,           Lucian%Durham%NULL%1,           Subbarao%Elapavaluru%NULL%1,           Carlos V.%Elzo Kraemer%NULL%1,           EDMÍLSON CARDOSO%Filho%NULL%1,           Jillian%Fitzgerald%NULL%1,           Giuseppe%Foti%NULL%1,           Matthew%Fox%NULL%1,           David%Fritschen%NULL%1,           David%Fullerton%NULL%1,           Elton%Gelandt%NULL%1,           Stacy%Gerle%NULL%1,           Marco%Giani%NULL%1,           Si Guim%Goh%NULL%1,           Sara%Govener%NULL%1,           Julie%Grone%NULL%1,           Miles%Guber%NULL%1,           Vadim%Gudzenko%NULL%1,           Daniel%Gutteridge%NULL%1,           Jennifer%Guy%NULL%1,           Jonathan%Haft%NULL%1,           Cameron%Hall%NULL%1,           Ibrahim Fawzy%Hassan%NULL%1,           Rubén%Herrán%NULL%1,           Hitoshi%Hirose%NULL%0,           Abdulsalam Saif%Ibrahim%NULL%1,           Don%Igielski%NULL%1,           Felicia A.%Ivascu%NULL%1,           Jaume%Izquierdo Blasco%NULL%1,           Julie%Jackson%NULL%1,           Harsh%Jain%NULL%1,           Bhavini%Jaiswal%NULL%1,           Andrea C.%Johnson%NULL%1,           Jenniver A.%Jurynec%NULL%1,           Norma M%Kellter%NULL%1,           Adam%Kohl%NULL%1,           Zachary%Kon%NULL%1,           Markus%Kredel%NULL%1,           Karen%Kriska%NULL%1,           Chandra%Kunavarapu%NULL%1,           Oude%Lansink-Hartgring%NULL%1,           Jeliene%LaRocque%NULL%1,           Sharon Beth%Larson%NULL%1,           Tracie%Layne%NULL%1,           Stephane%Ledot%NULL%1,           Napolitan%Lena%NULL%1,           Jonathan%Lillie%NULL%1,           Gösta%Lotz%NULL%1,           Mark%Lucas%NULL%1,           Lee%Ludwigson%NULL%1,           Jacinta J.%Maas%NULL%1,           Joanna%Maertens%NULL%1,           David%Mast%NULL%1,           Scott%McCardle%NULL%1,           Bernard%McDonald%NULL%1,           Allison%McLarty%NULL%1,           Chelsea%McMahon%NULL%1,           Patrick%Meybohm%NULL%1,           Bart%Meyns%NULL%1,           Casey%Miller%NULL%1,           Fernando%Moraes Neto%NULL%1,           Kelly%Morris%NULL%1,           Ralf%Muellenbach%NULL%1,           Meghan%Nicholson%NULL%1,           Serena%O''Brien%NULL%1,           Kathryn%O''Keefe%NULL%1,           Tawnya%Ogston%NULL%1,           Gary%Oldenburg%NULL%1,           Fabiana M.%Oliveira%NULL%1,           Emily%Oppel%NULL%1,           Diego%Pardo%NULL%2,           Diego%Pardo%NULL%0,           Sara J.%Parker%NULL%1,           Finn M.%Pedersen%NULL%1,           Crescens%Pellecchia%NULL%1,           Jose A.S.%Pelligrini%NULL%1,           Thao T.N.%Pham%NULL%1,           Ann R.%Phillips%NULL%1,           Tasneem%Pirani%NULL%1,           Paweł%Piwowarczyk%NULL%1,           Robert%Plambeck%NULL%1,           William%Pruett%NULL%1,           Brittany%Quandt%NULL%1,           Kollengode%Ramanathan%NULL%1,           Alejandro%Rey%NULL%1,           Christian%Reyher%NULL%1,           Jordi%Riera del Brio%NULL%1,           Rachel%Roberts%NULL%1,           David%Roe%NULL%1,           Peter P.%Roeleveld%NULL%1,           Janet%Rudy%NULL%1,           Luis F.%Rueda%NULL%1,           Emanuele%Russo%NULL%1,           Jesús%Sánchez Ballesteros%NULL%1,           Nancy%Satou%NULL%1,           Mauricio Guidi%Saueressig%NULL%1,           Paul C.%Saunders%NULL%1,           Margaret%Schlotterbeck%NULL%1,           Patricia%Schwarz%NULL%1,           Nicole%Scriven%NULL%1,           Alexis%Serra%NULL%1,           Mohammad%Shamsah%NULL%1,           Lucy%Sim%NULL%1,           Alexandra%Smart%NULL%1,           Adam%Smith%NULL%1,           Deane%Smith%NULL%1,           Maggie%Smith%NULL%1,           Neel%Sodha%NULL%1,           Michael%Sonntagbauer%NULL%1,           Marc%Sorenson%NULL%1,           Eric B%Stallkamp%NULL%1,           Allison%Stewart%NULL%1,           Kathy%Swartz%NULL%1,           Koji%Takeda%NULL%1,           Shaun%Thompson%NULL%1,           Bridget%Toy%NULL%1,           Divina%Tuazon%NULL%1,           Makoto%Uchiyama%NULL%1,           Obiora I.%Udeozo%NULL%1,           Scott%van Poppel%NULL%1,           Corey%Ventetuolo%NULL%1,           Leen%Vercaemst%NULL%1,           Nguyen V.%Vinh Chau%NULL%1,           I-Wen%Wang%NULL%1,           Carrie%Williamson%NULL%1,           Brock%Wilson%NULL%1,           Helen%Winkels%NULL%1]'
$ws.Range("E5").Value = '[Matthew%Charlton%NULL%1,           Susan%Dashey%NULL%1,           Abigail%Stubbs%NULL%1,           Florence YL%Lai%NULL%1,           Paul W%Bird%NULL%1,           Vinay%Badhwar%NULL%1,           Julian W%Tang%NULL%1]'
$ws.Range("E6").Value = '[Nicolas%Cousin%Nicolas.COUSIN@CHRU-LILLE.FR%1,           Claire%Bourel%claire.bourel.etu@univ-lille.fr%1,           Dorothee%Carpentier%Dorothee.Carpentier@chu-rouen.fr%1,           Julien%Goutay%Julien.LABREUCHE@CHRU-LILLE.FR%1,           Agnes%Mugnier%Agnes.MUGNIER@CHRU-LILLE.FR%1,           Julien%Labreuche%Julien.LABREUCHE@CHRU-LILLE.FR%1,           Elise%Godeau%elisegodeau@gmail.com%1,           Thomas%Clavier%thomas.clavier@chu-rouen.fr%1,           Steven%Grange%Steven.Grange@chu-rouen.fr%1,           Fabienne%Tamion%Fabienne.Tamion@chu-rouen.fr%1,           Arthur%Durand%Arthur.DURAND@CHRU-LILLE.FR%1,           Mouhamed D.%Moussa%Mouhamed.MOUSSA@CHRU-LILLE.FR%1,           Thibault%Duburcq%NULL%1]'
$ws.Range("E7").Value = '[Pierre-Emmanuel%Falcoz%NULL%1,           Alexandra%Monnier%NULL%1,           Marc%Puyraveau%NULL%1,           Stephanie%Perrier%NULL%1,           Pierre-Olivier%Ludes%NULL%1,           Anne%Olland%NULL%1,           Paul-Michel%Mertes%NULL%2,           Francis%Schneider%NULL%2,           Julie%Helms%NULL%3,           NULL%NULL%NULL%0]'
$ws.Range("E8").Value = '[Julien%Guihaire%xref no email%1,    Clark g.%Owyang%xref no email%1,    Jai%Madhok%xref no email%1,    Florent%Laverdure%xref no email%1,    Ma\u00efra%Gaillard%xref no email%1,    Antoine%Girault%xref no email%1,    Guillaume%Lebreton%xref no email%1,    Olaf%Mercier%xref no email%1]'
$ws.Range("E9").Value = '[Pierre%Huette%huette.pierre@chu-amiens.fr%1,           Christophe%Beyls%NULL%2,           Christophe%Beyls%NULL%0,           Mathieu%Guilbart%NULL%1,           Alexandre%Coquet%NULL%1,           Pascal%Berna%NULL%1,           Guillaume%Haye%NULL%1,           Pierre-Alexandre%Roger%NULL%1,           Patricia%Besserve%NULL%1,           Michael%Bernasinski%NULL%1,           Hervé%Dupont%NULL%1,           Osama%Abou-Arab%NULL%1,           Yazine%Mahjoub%NULL%1]'
$ws.Range("E10").Value = '[Markus%Jäckel%markus.jaeckel@universitaets-herzzentrum.de%1,           Jonathan%Rilinger%NULL%2,           Jonathan%Rilinger%NULL%0,           Corinna Nadine%Lang%NULL%2,           Corinna Nadine%Lang%NULL%0,           Viviane%Zotzmann%NULL%1,           Klaus%Kaier%NULL%1,           Peter%Stachon%NULL%1,           Paul Marc%Biever%NULL%1,           Tobias%Wengenmayer%NULL%1,           Daniel%Duerschmied%NULL%1,           Christoph%Bode%NULL%1,           Dawid Leander%Staudacher%NULL%1,           Alexander%Supady%NULL%1]'
$ws.Range("E11").Value = '[Woo Sung%Jang%NULL%1,           JaeBum%Kim%NULL%1,           Jonghyun%Baek%NULL%1,           Hanna%Jung%NULL%1,           Jae Seok%Jang%NULL%1,           Jae Seok%Park%NULL%0,           Tak-Hyuk%Oh%NULL%1,           Se Yong%Jang%NULL%1,           Yun Seok%Kim%NULL%1,           Yong Shik%Kwon%NULL%0]'
$ws.Range("E12").Value = '[C.%Le Breton%NULL%1,           S.%Besset%NULL%1,           S.%Freita-Ramos%NULL%1,           M.%Amouretti%NULL%1,           P.A.%Billiet%NULL%1,           M.%Dao%NULL%1,           L.M.%Dumont%NULL%1,           L.%Federici%NULL%1,           B.%Gaborieau%NULL%1,           D.%Longrois%NULL%1,           P.%Postel-Vinay%NULL%1,           C.%Vuillard%NULL%1,           N.%Zucman%NULL%1,           G.%Lebreton%NULL%1,           A.%Combes%NULL%1,           D.%Dreyfuss%NULL%1,           J.D.%Ricard%NULL%1,           D.%Roux%NULL%1]'
$ws.Range("E13").Value = '[Mathieu%Jozwiak%NULL%1,           Jean-Daniel%Chiche%NULL%1,           Julien%Charpentier%NULL%1,           Zakaria%Ait Hamou%NULL%1,           Paul%Jaubert%NULL%1,           Sarah%Benghanem%NULL%1,           Pierre%Dupland%NULL%1,           Ariane%Gavaud%NULL%1,           Frédéric%Péne%NULL%1,           Alain%Cariou%NULL%1,           Jean-Paul%Mira%NULL%1,           Lee S.%Nguyen%NULL%1]'
$ws.Range("E14").Value = '[J.%Masur%xref no email%1,    C.W.%Freeman%xref no email%1,    S.%Mohan%xref no email%1]'
$ws.Range("E15").Value = '[Kevin%Roedl%NULL%0,           Dominik%Jarczak%NULL%1,           Liina%Thasler%NULL%1,           Martin%Bachmann%NULL%1,           Frank%Schulte%NULL%1,           Berthold%Bein%NULL%1,           Christian Friedrich%Weber%NULL%1,           Ulrich%Schäfer%NULL%1,           Carsten%Veit%NULL%1,           Hans-Peter%Hauber%NULL%1,           Sebastian%Kopp%NULL%1,           Karsten%Sydow%NULL%1,           Andreas%de Weerth%NULL%1,           Marc%Bota%NULL%1,           Rüdiger%Schreiber%NULL%1,           Oliver%Detsch%NULL%1,           Jan-Peer%Rogmann%NULL%1,           Daniel%Frings%NULL%1,           Barbara%Sensen%NULL%1,           Christoph%Burdelski%NULL%1,           Olaf%Boenisch%NULL%1,           Axel%Nierhaus%NULL%1,           Geraldine%de Heer%NULL%1,           Stefan%Kluge%NULL%1]'
$ws.Range("E16").Value = '[Emily%Shih%NULL%1,           J. Michael%DiMaio%NULL%1,           John J.%Squiers%NULL%1,           Jasjit K.%Banwait%NULL%1,           Dan M.%Meyer%NULL%1,           Timothy J.%George%NULL%1,           Gary S.%Schwartz%NULL%1,           Britton A.%Blough%NULL%1,           Robert L.%Gottlieb%NULL%1,           Gonzo%Gonzalez-Stawinski%NULL%1,           Anita R.%Krueger%NULL%1,           Kara%Monday%NULL%1,           Ramachandra C.%Reddy%NULL%1,           Jorge F.%Velazco%NULL%1,           Omar%Hernandez%NULL%1,           Chris%Martin%NULL%1,           Jenelle%Sheasby%NULL%1]'
$ws.Range("E17").Value = '[Joe%Zhang%NULL%1,           Blair%Merrick%NULL%2,           Blair%Merrick%NULL%0,           Genex L.%Correa%NULL%1,           Luigi%Camporota%NULL%1,           Andrew%Retter%NULL%1,           Andrew%Doyle%NULL%1,           Guy W.%Glover%NULL%2,           Guy W.%Glover%NULL%0,           Peter B.%Sherren%NULL%1,           Stephen J.%Tricklebank%NULL%1,           Sangita%Agarwal%NULL%1,           Boris E.%Lams%NULL%1,           Nicholas A.%Barrett%NULL%1,           Nicholas%Ioannou%NULL%1,           Jonathan%Edgeworth%NULL%1,           Christopher I.S.%Meadows%NULL%1]'
